{"js": "// Merge the three runs that read:\n//   \"Drag and drop the component into your frame\" + \" while holding Alt on windows\" + \".\"\n// into a single run that reads:\n//   \"Drag and drop the component into your frame.\"\n// (i.e. drop \" while holding Alt on windows\" and keep a single trailing period).\n\nconst body = context.document.body;\n\nconst oldText = \"Drag and drop the component into your frame while holding Alt on windows.\";\nconst newText = \"Drag and drop the component into your frame.\";\n\nconst results = body.search(oldText, { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target sentence to edit: \" + oldText);\n}\n\n// Replacing the whole matched range collapses it back down to the first\n// run's formatting (rFonts/color/sz/szCs), which is what every run in this\n// sentence already shared, and removes the now-empty trailing runs.\nresults.items[0].insertText(newText, \"Replace\");\n\nawait context.sync();\n", "ps1": "# Remove \" while holding Alt on windows\" from the sentence:\n#   \"Drag and drop the component into your frame while holding Alt on windows.\"\n# leaving it as:\n#   \"Drag and drop the component into your frame.\"\n#\n# Locating the exact phrase with Find and deleting just that sub-range (rather\n# than rewriting the whole sentence) keeps the surrounding runs' formatting\n# (and the first run's rsid) untouched, and Word naturally coalesces the\n# now-adjacent, identically-formatted runs on either side of the deletion.\n\n$d = $word.ActiveDocument\n\n$r = $d.Range(0, $d.Content.End)\n$found = $r.Find.Execute(\" while holding Alt on windows\")\n\nif (-not $found) {\n    throw \"Could not find target phrase ' while holding Alt on windows' to remove.\"\n}\n\n$r.Delete()\n"}
